$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.459.13'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.573.04'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.04'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3738'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.97'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3402'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07553'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.142'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.991'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.944'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.566.24'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001122'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06728'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.252'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.42'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.465.12'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.352'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.583'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -6.13%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.67'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.020'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.85'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.744.63'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.130'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.53%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.830'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08414'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.375'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02461'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2294'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.471'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.36'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6274'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.92'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.809'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5845'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.15%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.87'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.34%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.086'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.227'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07328'
